{"js": "// Replace the division-problem text in the worksheet table cells.\n// Each pair is the exact original cell text (unique in the document)\n// and its replacement, taken from the OOXML diff.\nconst replacements = [\n  [\"640\u00f72=\", \"333\u00f76=\"],\n  [\"920\u00f76=\", \"290\u00f76=\"],\n  [\"598\u00f72=\", \"238\u00f73=\"],\n  [\"581\u00f79=\", \"944\u00f77=\"],\n  [\"460\u00f75=\", \"375\u00f77=\"],\n  [\"803\u00f77=\", \"830\u00f72=\"],\n  [\"507\u00f77=\", \"712\u00f76=\"],\n  [\"586\u00f72=\", \"238\u00f73=\"],\n  [\"236\u00f79=\", \"535\u00f76=\"],\n  [\"863\u00f78=\", \"463\u00f77=\"],\n  [\"874\u00f77=\", \"579\u00f77=\"],\n  [\"258\u00f76=\", \"166\u00f74=\"],\n  [\"523\u00f72=\", \"653\u00f79=\"],\n  [\"228\u00f78=\", \"303\u00f78=\"],\n  [\"905\u00f75=\", \"378\u00f73=\"],\n  [\"296\u00f75=\", \"842\u00f75=\"],\n  [\"832\u00f75=\", \"348\u00f78=\"],\n  [\"154\u00f72=\", \"437\u00f74=\"],\n  [\"729\u00f78=\", \"774\u00f73=\"],\n  [\"518\u00f72=\", \"818\u00f79=\"],\n  [\"405\u00f79=\", \"306\u00f76=\"],\n  [\"960\u00f76=\", \"206\u00f77=\"],\n  [\"284\u00f73=\", \"331\u00f72=\"],\n  [\"800\u00f75=\", \"378\u00f79=\"],\n  [\"284\u00f79=\", \"844\u00f78=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the division-problem text in the worksheet table cells.\n# Each pair is the exact original cell text (unique in the document)\n# and its replacement, taken from the OOXML diff.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"640\u00f72=\", \"333\u00f76=\"),\n    @(\"920\u00f76=\", \"290\u00f76=\"),\n    @(\"598\u00f72=\", \"238\u00f73=\"),\n    @(\"581\u00f79=\", \"944\u00f77=\"),\n    @(\"460\u00f75=\", \"375\u00f77=\"),\n    @(\"803\u00f77=\", \"830\u00f72=\"),\n    @(\"507\u00f77=\", \"712\u00f76=\"),\n    @(\"586\u00f72=\", \"238\u00f73=\"),\n    @(\"236\u00f79=\", \"535\u00f76=\"),\n    @(\"863\u00f78=\", \"463\u00f77=\"),\n    @(\"874\u00f77=\", \"579\u00f77=\"),\n    @(\"258\u00f76=\", \"166\u00f74=\"),\n    @(\"523\u00f72=\", \"653\u00f79=\"),\n    @(\"228\u00f78=\", \"303\u00f78=\"),\n    @(\"905\u00f75=\", \"378\u00f73=\"),\n    @(\"296\u00f75=\", \"842\u00f75=\"),\n    @(\"832\u00f75=\", \"348\u00f78=\"),\n    @(\"154\u00f72=\", \"437\u00f74=\"),\n    @(\"729\u00f78=\", \"774\u00f73=\"),\n    @(\"518\u00f72=\", \"818\u00f79=\"),\n    @(\"405\u00f79=\", \"306\u00f76=\"),\n    @(\"960\u00f76=\", \"206\u00f77=\"),\n    @(\"284\u00f73=\", \"331\u00f72=\"),\n    @(\"800\u00f75=\", \"378\u00f79=\"),\n    @(\"284\u00f79=\", \"844\u00f78=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
